$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the existing data rows (below the header) so we can rewrite
# them in the new order/content, including the newly added rows.
$ws.Range("A2:F20").ClearContents()

# New data, reordered by id (skill-score based ranking) and with new
# candidates/projects appended (RAPID, Alcheringa Pass Portal,
# TEDxIITGuwahati Website, Video Conferencing Project).
$data = @(
    @(1, "LLMGuard", 81.62, 1, 81.62, 2),
    @(1, "Literature Society IITJ Website", 78.28, 0.9, 70.45, 2),
    @(2, "SMART SENSING MIDDLEWARE", 110.08, 1, 100, 2),
    @(2, "RAPID", 81.62, 1, 81.62, 2),
    @(3, "LLMGuard", 85.02, 1, 85.02, 2),
    @(3, "Website for the Literature Society of the college", 78.28, 0.9, 70.45, 2),
    @(4, "Alcheringa Pass Portal", 85.02, 1, 85.02, 3),
    @(4, "TEDxIITGuwahati Website", 78.28, 0.85, 66.54000000000001, 3),
    @(4, "Video Conferencing Project", 78.28, 0.85, 66.54000000000001, 3)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
